$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed in the repulled data.
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -7
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = -1
$ws.Range("F6").Value = 13
$ws.Range("F8").Value = -1
$ws.Range("F10").Value = 1
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = -5
$ws.Range("F14").Value = -1
$ws.Range("F15").Value = -3
$ws.Range("F16").Value = 1
